$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Emner")

# --- Column A: "Id" -> "EmneNr", numeric topic ids -> "Topic N" text labels ---
$ws.Range("A1").Value = "EmneNr"

$ids = @(0,1,3,4,5,7,8,9,10,11,12,13,14,15,16,17,18,19)
for ($i = 0; $i -lt $ids.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = ("Topic " + $ids[$i])
}

# --- Column D: stamp every data row with the "SidstOpdatret" date ---
# Rows 2-16 already carry the date-formatted style (s="10"); rows 17-19 still
# use the older unformatted style (s="9"), so pull the date format from row 2
# onto them first, then write the date value into the whole D2:D19 range.
$ws.Range("D2").Copy()
$ws.Range("D17:D19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

for ($r = 2; $r -le 19; $r++) {
    $ws.Cells.Item($r, 4).Value = 45863
}

# --- Column A width: drop bestFit, widen slightly ---
$ws.Columns("A").ColumnWidth = 12.92

# --- Selection moves from C1 to B2 ---
[void]$ws.Range("B2").Select()
